# ST3 Line Rejection sheet update
#  - Row 5, columns C:G: re-type values from text ("234") to real numbers (234)
#  - Append new data rows 6 & 7 (SHIFT2 @ 2025-02-14T18:45) and row 8 (SHIFT2 @ 2025-02-06T18:47)
#  - Sheet used range grows from A1:J5 to A1:J8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must be stored as TEXT even when it looks like a
# number (e.g. "123"), mirroring how Excel keeps a cell textual when it was
# formatted as Text before typing. Re-applying the "Normal" style afterwards
# drops the temporary Text number-format so no stray formatting is left
# behind on the cell.
function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# --- Row 5: re-type C5:G5 as numbers (currently stored as text "234") ---
$ws.Cells.Item(5, 3).Value = 234
$ws.Cells.Item(5, 4).Value = 234
$ws.Cells.Item(5, 5).Value = 234
$ws.Cells.Item(5, 6).Value = 234
$ws.Cells.Item(5, 7).Value = 234

# --- Row 6 (new) ---
Set-TextValue 6 1 "2025-02-14T18:45"
Set-TextValue 6 2 "SHIFT2"
$ws.Cells.Item(6, 3).Value = 123
$ws.Cells.Item(6, 4).Value = 23
$ws.Cells.Item(6, 5).Value = 123
$ws.Cells.Item(6, 6).Value = 123
$ws.Cells.Item(6, 7).Value = 123
Set-TextValue 6 8 "123"
Set-TextValue 6 9 "123"
Set-TextValue 6 10 "1123"

# --- Row 7 (new) ---
Set-TextValue 7 1 "2025-02-14T18:45"
Set-TextValue 7 2 "SHIFT2"
$ws.Cells.Item(7, 3).Value = 123
$ws.Cells.Item(7, 4).Value = 23
$ws.Cells.Item(7, 5).Value = 123
$ws.Cells.Item(7, 6).Value = 123
$ws.Cells.Item(7, 7).Value = 123
Set-TextValue 7 8 "123"
Set-TextValue 7 9 "123"
Set-TextValue 7 10 "1123"

# --- Row 8 (new) ---
Set-TextValue 8 1 "2025-02-06T18:47"
Set-TextValue 8 2 "SHIFT2"
Set-TextValue 8 3 "12"
Set-TextValue 8 4 "1"
Set-TextValue 8 5 "2"
Set-TextValue 8 6 "1"
Set-TextValue 8 7 "2"
Set-TextValue 8 8 "aEF"
Set-TextValue 8 9 "SDG"
Set-TextValue 8 10 "ASDF"
